# LOB1201.docx edit: rotate several blocks of body text between paragraph
# slots (Objetivos / Docente(s) / Programa resumido / Programa / Avaliação /
# Bibliografia) while leaving every paragraph's style & run formatting
# (bold labels, italics, bullet lists, line breaks) untouched.
#
# Strategy: for every paragraph whose text changes, grab a duplicate of
# that paragraph's Range and do a scoped Find & Replace (Wrap = wdFindStop
# = 0, so the search can never leak into a neighbouring paragraph). A
# vertical-tab character (chr 11) inside the search/replacement text
# matches/produces a <w:br/> line break, which lets us collapse or expand
# multi-line runs exactly like Word's own Find/Replace does.

$d = $word.ActiveDocument
$vt = [char]11

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$Old,
        [string]$New
    )
    $para = $d.Paragraphs.Item($Index)
    $rng = $para.Range.Duplicate
    $ok = $rng.Find.Execute($Old, $true, $false, $false, $false, $false, $true, 0, $false, $New, 2)
    if (-not $ok) {
        Write-Output ("WARNING: replacement not found in paragraph " + $Index + ": " + $Old)
    }
}

# --- Paragraph 6: "Objetivos" body (PT) ---------------------------------
Replace-InParagraph 6 `
    "Dar conhecimentos aos alunos de noções básicas sobre ecologia e impacto das atividades da engenharia sobre o meio ambiente. Conceitos legais e institucionais para o desenvolvimento sustentável." `
    "Conceitos e Definições. Questões Ambientais. Desenvolvimento Sustentável. Desempenho Ambiental. Processos Ambientais. Norma Ambiental."

# --- Paragraph 7: "Objetivos" body (EN, italic) -------------------------
Replace-InParagraph 7 `
    "Give students knowledge of the basics of ecology and impact of engineering activities on the environment. Legal and institutional concepts for sustainable development." `
    "Concepts and Definitions. Environmental issues. Sustainable Development. Environmental performance. Environmental processes. Environmental standard."

# --- Paragraph 9: "Docente(s) Responsável(eis)" bullet list ------------
Replace-InParagraph 9 `
    "9146830 - Danúbia Caporusso Bargos" `
    "Dar conhecimentos aos alunos de noções básicas sobre ecologia e impacto das atividades da engenharia sobre o meio ambiente. Conceitos legais e institucionais para o desenvolvimento sustentável."

Replace-InParagraph 9 `
    "5464150 - Mariana Consiglio Kasemodel" `
    "CONCEITOS E DEFINIÇÕES. Engenharia Ambiental. Meio Ambiente. Poluição Ambiental. Componentes Ambientais Críticos. QUESTÕES AMBIENTAIS. O Sujeito das Transformações Ambientais. Energia e o Meio Ambiente. Impactos Ambientais nos Três Meios. Equilíbrio Ameaçado. DESENVOLVIMENTO SUSTENTÁVEL. Conceitos Básicos. Aspectos legais. DESEMPENHO AMBIENTAL. Monitoramento Ambiental. Abrangência do Desenvolvimento Sustentável. Definição de Indicadores. Definição de Indicadores Sustentáveis. Indicadores de Desenvolvimento Humano – IDH. Indicadores de Sustentabilidade Ambiental. Controle de Processos Ambientais. PROCESSOS AMBIENTAIS. Controle Processo ETA. Água na Natureza. Caracterização da Água. Indicadores de Qualidade da Água. NORMA AMBIENTAL. Portaria 518. CONAMA 20. Desastre Ecológico. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina."

# --- Paragraph 11: "Programa resumido" body (PT) ------------------------
Replace-InParagraph 11 `
    "Conceitos e Definições. Questões Ambientais. Desenvolvimento Sustentável. Desempenho Ambiental. Processos Ambientais. Norma Ambiental." `
    "Aulas expositivas com a utilização de recursos de projeções e audiovisual."

# --- Paragraph 12: "Programa resumido" body (EN, italic) ---------------
Replace-InParagraph 12 `
    "Concepts and Definitions. Environmental issues. Sustainable Development. Environmental performance. Environmental processes. Environmental standard." `
    "Give students knowledge of the basics of ecology and impact of engineering activities on the environment. Legal and institutional concepts for sustainable development."

# --- Paragraph 14: "Programa" body (PT) ---------------------------------
Replace-InParagraph 14 `
    "CONCEITOS E DEFINIÇÕES. Engenharia Ambiental. Meio Ambiente. Poluição Ambiental. Componentes Ambientais Críticos. QUESTÕES AMBIENTAIS. O Sujeito das Transformações Ambientais. Energia e o Meio Ambiente. Impactos Ambientais nos Três Meios. Equilíbrio Ameaçado. DESENVOLVIMENTO SUSTENTÁVEL. Conceitos Básicos. Aspectos legais. DESEMPENHO AMBIENTAL. Monitoramento Ambiental. Abrangência do Desenvolvimento Sustentável. Definição de Indicadores. Definição de Indicadores Sustentáveis. Indicadores de Desenvolvimento Humano – IDH. Indicadores de Sustentabilidade Ambiental. Controle de Processos Ambientais. PROCESSOS AMBIENTAIS. Controle Processo ETA. Água na Natureza. Caracterização da Água. Indicadores de Qualidade da Água. NORMA AMBIENTAL. Portaria 518. CONAMA 20. Desastre Ecológico. A disciplina pode contar com viagens didáticas para complementação do conteúdo da disciplina." `
    "Média ponderada de 2 avaliações escritas com nota final (NF ≥ 5,0)"

# --- Paragraph 17: "Avaliação" bullet list ------------------------------
# Run "Método: " value
Replace-InParagraph 17 `
    "Aulas expositivas com a utilização de recursos de projeções e audiovisual." `
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova escrita de recuperação a ser aplicada"

# Run "Critério: " value -> expands into the 5-item bibliography list
$biblio1 = "1)        BRAGA, B.; HESPANHOL, I.; CONEJO, J. G. L.; MIERZWA, J. C.; BARROS, M. T. L.; SPENCER, M.; PORTO, M.; NUCCI, N.; JULIANO, N.; EIGER, S. Introdução à Engenharia Ambiental: O Desafio do Desenvolvimento Sustentável. Pearson (2ª Edição), 336 p., 2005."
$biblio2 = "2)        VESILIND, P.A.; MORGAN, S. M.; HEINE, L. G. Introdução à Engenharia Ambiental. Cengage (3ª edição), 472 p., 2018."
$biblio3 = "3)        CALIJURI, M. C.; CUNHA, D. G. F. Engenharia Ambiental: Conceitos, Tecnologias e Gestão. Elsevier (1ª Edição), 832 p., 2012."
$biblio4 = "4)        CAPAZ, R. S.; HORTA NOGUEIRA, L. A. Ciências Ambientais para Engenharia. Elsevier (1ª Edição), 252 p., 2014."
$biblio5 = "5)        DAVIS, M. L.; MASTEN, S. J. Princípios de Engenharia Ambiental. Mc Graw Hill Educations (3ª Edição), 872 p., 2016;"
$biblioAll = $biblio1 + $vt + $biblio2 + $vt + $biblio3 + $vt + $biblio4 + $vt + $biblio5

Replace-InParagraph 17 `
    "Média ponderada de 2 avaliações escritas com nota final (NF ≥ 5,0)" `
    $biblioAll

# Run "Norma de recuperação: " value
Replace-InParagraph 17 `
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova escrita de recuperação a ser aplicada" `
    "9146830 - Danúbia Caporusso Bargos"

# --- Paragraph 19: "Bibliografia" body ----------------------------------
# Collapses the 5-item list (4 internal line breaks) back down to a single
# teacher-name line.
$oldBiblioAll = $biblio1 + $vt + $biblio2 + $vt + $biblio3 + $vt + $biblio4 + $vt + $biblio5
Replace-InParagraph 19 `
    $oldBiblioAll `
    "5464150 - Mariana Consiglio Kasemodel"

Write-Output "Done."
